$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B values for rows 59-66 (alternating 0.015 / 0.158 pattern)
$ws.Range("B59").Value = 0.015
$ws.Range("B60").Value = 0.158
$ws.Range("B61").Value = 0.015
$ws.Range("B62").Value = 0.158
$ws.Range("B63").Value = 0.015
$ws.Range("B64").Value = 0.158
$ws.Range("B65").Value = 0.015
$ws.Range("B66").Value = 0.158

# Delete rows 67-82 (trailing duplicate param block that's no longer needed)
$ws.Range("A67:Q82").EntireRow.Delete()

$wb.Save()
